# Update existing rows 215-232 (columns D, K, L, M, N, O, P, Q, S) and append
# two brand-new rows (233, 234) at the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: copy the date number-format from an existing date cell so new
# rows keep the same "YYYY-MM-DD HH:MM:SS" style used by column D. ---
$dateFormat = $ws.Cells.Item(2, 4).NumberFormat

# -----------------------------------------------------------------------
# Updated values for existing rows 215-232.
# Columns: D=4 (Fecha), K=11 (Variedad), L=12 (Calidad), M=13 (Volumen),
#          N=14 (Precio minimo), O=15 (Precio maximo),
#          P=16 (Precio promedio ponderado), Q=17 (Unidad de comercializacion),
#          S=19 (Precio $/Kg)
# -----------------------------------------------------------------------
$updates = @(
    @{Row=215; D=44461; K='Navel Late'; L='Primera'; M=240; N=6000;  O=6500;  P=6250;  Q='$/bandeja 15 kilos granel'; S=417},
    @{Row=216; D=44461; K='Navel Late'; L='Segunda'; M=180; N=5000;  O=5500;  P=5250;  Q='$/bandeja 15 kilos granel'; S=350},
    @{Row=217; D=44357; K='Fukumoto';   L='Primera'; M=300; N=11000; O=12000; P=11500; Q='$/bandeja 15 kilos granel'; S=767},
    @{Row=218; D=44162; K='Valencia';   L='Primera'; M=120; N=14000; O=15000; P=14500; Q='$/caja 15 kilos granel';    S=967},
    @{Row=219; D=44410; K='Navel Late'; L='Primera'; M=400; N=6000;  O=6500;  P=6250;  Q='$/bandeja 15 kilos granel'; S=417},
    @{Row=220; D=44410; K='Navel Late'; L='Segunda'; M=300; N=5000;  O=5500;  P=5250;  Q='$/bandeja 15 kilos granel'; S=350},
    @{Row=221; D=44411; K='Navel Late'; L='Primera'; M=400; N=6500;  O=7000;  P=6750;  Q='$/bandeja 15 kilos granel'; S=450},
    @{Row=222; D=44411; K='Navel Late'; L='Segunda'; M=80;  N=6000;  O=6000;  P=6000;  Q='$/bandeja 15 kilos granel'; S=400},
    @{Row=223; D=44176; K='Valencia';   L='Primera'; M=120; N=13000; O=14000; P=13417; Q='$/bandeja 15 kilos granel'; S=894},
    @{Row=224; D=44376; K='Fukumoto';   L='Primera'; M=120; N=8000;  O=8500;  P=8250;  Q='$/bandeja 15 kilos granel'; S=550},
    @{Row=225; D=44376; K='Fukumoto';   L='Segunda'; M=120; N=7000;  O=7500;  P=7250;  Q='$/bandeja 15 kilos granel'; S=483},
    @{Row=226; D=44292; K='Valencia';   L='Primera'; M=40;  N=17500; O=18000; P=17750; Q='$/bandeja 15 kilos granel'; S=1183},
    @{Row=227; D=44358; K='Fukumoto';   L='Primera'; M=120; N=10000; O=11000; P=10500; Q='$/caja 15 kilos granel';    S=700},
    @{Row=228; D=44211; K='Valencia';   L='Primera'; M=100; N=18000; O=19000; P=18450; Q='$/bandeja 15 kilos granel'; S=1230},
    @{Row=229; D=44425; K='Fukumoto';   L='Primera'; M=200; N=6000;  O=6200;  P=6100;  Q='$/bandeja 15 kilos granel'; S=407},
    @{Row=230; D=44425; K='Fukumoto';   L='Segunda'; M=120; N=5000;  O=5500;  P=5250;  Q='$/bandeja 15 kilos granel'; S=350},
    @{Row=231; D=44425; K='Navel Late'; L='Primera'; M=200; N=5000;  O=5500;  P=5250;  Q='$/bandeja 15 kilos granel'; S=350},
    @{Row=232; D=44425; K='Navel Late'; L='Segunda'; M=100; N=4500;  O=4800;  P=4650;  Q='$/bandeja 15 kilos granel'; S=310}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D    # D - Fecha
    $ws.Cells.Item($r, 11).Value = $u.K    # K - Variedad
    $ws.Cells.Item($r, 12).Value = $u.L    # L - Calidad
    $ws.Cells.Item($r, 13).Value = $u.M    # M - Volumen
    $ws.Cells.Item($r, 14).Value = $u.N    # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $u.O    # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $u.P    # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $u.Q    # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $u.S    # S - Precio $/Kg
}

# -----------------------------------------------------------------------
# Two brand-new rows appended at the end of the table (233 and 234), with
# the constant columns (A, B, C, E, F, G, H, I, J, R, T) carried over from
# the surrounding rows.
# -----------------------------------------------------------------------
$newRows = @(
    @{Row=233; D=44323; K='Fukumoto'; L='Primera'; M=120; N=13000; O=14000; P=13500; Q='$/bandeja 15 kilos granel'; S=900},
    @{Row=234; D=44323; K='Fukumoto'; L='Segunda'; M=80;  N=12000; O=12000; P=12000; Q='$/bandeja 15 kilos granel'; S=800}
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value  = 7
    $ws.Cells.Item($r, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
    $ws.Cells.Item($r, 3).Value  = 'Ñuble'
    $ws.Cells.Item($r, 4).Value  = $nr.D
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 5).Value  = 16
    $ws.Cells.Item($r, 6).Value  = 'Fruta'
    $ws.Cells.Item($r, 7).Value  = 100102
    $ws.Cells.Item($r, 8).Value  = 'Cítricos'
    $ws.Cells.Item($r, 9).Value  = 100102005
    $ws.Cells.Item($r, 10).Value = 'Naranja'
    $ws.Cells.Item($r, 11).Value = $nr.K
    $ws.Cells.Item($r, 12).Value = $nr.L
    $ws.Cells.Item($r, 13).Value = $nr.M
    $ws.Cells.Item($r, 14).Value = $nr.N
    $ws.Cells.Item($r, 15).Value = $nr.O
    $ws.Cells.Item($r, 16).Value = $nr.P
    $ws.Cells.Item($r, 17).Value = $nr.Q
    $ws.Cells.Item($r, 18).Value = 'Región de O''Higgins'
    $ws.Cells.Item($r, 19).Value = $nr.S
    $ws.Cells.Item($r, 20).Value = 15
}
